$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.004.21'
$ws.Range("D3").Value = '1.827.72'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '311.62'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").Value = '1.007'
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").Value = '0.4661'
$ws.Range("E7").Value = '  -0.97%  '
$ws.Range("D8").Value = '0.3714'
$ws.Range("E8").Value = '  +1.67%  '
$ws.Range("D9").Value = '0.07362'
$ws.Range("D10").Value = '0.8739'
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("D11").Value = '19.96'
$ws.Range("E11").Value = '  -1.87%  '
$ws.Range("D12").Value = '0.07826'
$ws.Range("E12").Value = '  +6.91%  '
$ws.Range("D13").Value = '1.842.01'
$ws.Range("E13").Value = '  -4.39%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.357'
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '6.564'
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("E16").Value = '  -1.47%  '
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("D18").Value = '0.000008888'
$ws.Range("E18").Value = '  +2.06%  '
$ws.Range("D19").Value = '1.007'
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("D21").Value = '26.835.49'
$ws.Range("E21").Value = '  -3.20%  '
$ws.Range("D22").Value = '5.155'
$ws.Range("E22").Value = '  -1.61%  '
$ws.Range("D23").Value = '10.58'
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").Value = '2.074.92'
$ws.Range("E24").Value = '  -1.09%  '
$ws.Range("D25").Value = '152.59'
$ws.Range("E25").Value = '  +0.60%  '
$ws.Range("D26").Value = '1.834'
$ws.Range("D27").Value = '18.27'
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("E28").Value = '  -1.96%  '
$ws.Range("D29").Value = '5.125'
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("D30").Value = '115.46'
$ws.Range("E30").Value = '  -0.50%  '
$ws.Range("D31").Value = '0.08873'
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").Value = '2.975'
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("D33").Value = '0.7294'
$ws.Range("E33").Value = '  -1.64%  '
$ws.Range("D34").Value = '4.444'
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("E35").Value = '  -2.48%  '
$ws.Range("D36").Value = '2.525'
$ws.Range("E36").Value = '  +4.96%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.01958'
$ws.Range("E37").Value = '  +0.63%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '1.076'
$ws.Range("E38").Value = '  -1.16%  '
$ws.Range("D39").Value = '0.05239'
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("D40").Value = '7.290'
$ws.Range("E40").Value = '  +0.98%  '
$ws.Range("D41").Value = '2.927'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = '0.5196'
$ws.Range("E42").Value = '  -0.98%  '
$ws.Range("D43").Value = '0.8658'
$ws.Range("E43").Value = '  -14.24%  '
$ws.Range("D44").Value = '0.1626'
$ws.Range("E44").Value = '  -0.93%  '
$ws.Range("D45").Value = '8.231'
$ws.Range("D46").Value = '0.4851'
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '10.25'
$ws.Range("E47").Value = '  -0.87%  '
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '1.008'
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("D49").Value = '102.71'
$ws.Range("E49").Value = '  -1.37%  '
$ws.Range("D50").Value = '1.625'
$ws.Range("E50").Value = '  -1.66%  '
$ws.Range("E51").Value = '  -1.30%  '
